$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# New file handed off for localization: a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md
# This takes the row previously used by ".localization-config" on every
# sheet, and ".localization-config" moves down to a new row 4.
# -----------------------------------------------------------------------

$mdOld    = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$mdNew    = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
$config   = ".localization-config"
$ready    = "Ready for handoff"
$notLoc   = "Not to be localized"
$include  = "Include"
$ignored  = "Ignored"
$epoch    = "0001-01-01 00:00:00"

$zhcnOldXlf = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf"
$zhcnNewXlf = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf"
$dedeOldXlf = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf"
$dedeNewXlf = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf"

$zhcnNewDt = "2016-02-22 13:47:00"
$dedeNewDt = "2016-02-22 13:47:15"

$mdUrlOld = "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/e2e/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$mdUrlNew = "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/e2e/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/.localization-config"

$zhcnXlfUrlOld = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e30bce61b2ccf4fd628866be3658b41fd2c9f242/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf"
$zhcnXlfUrlNew = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e30bce61b2ccf4fd628866be3658b41fd2c9f242/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf"

$dedeXlfUrlOld = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b61d5a74457a7aa614eb0ad167fdc561f96b011/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf"
$dedeXlfUrlNew = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b61d5a74457a7aa614eb0ad167fdc561f96b011/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf"

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = $mdNew
$ws.Range("B3").Value = $ready
$ws.Range("C3").Value = $ready

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = $notLoc
$ws.Range("C4").Value = $notLoc

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlOld, "", "", $mdOld)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlNew, "", "", $mdNew)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $config)

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = $mdNew
$ws.Range("B3").Value = $ready
$ws.Range("C3").Value = $zhcnNewXlf
$ws.Range("D3").Value = $zhcnNewDt
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $include

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = $notLoc
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = $ignored

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlOld, "", "", $mdOld)
$ws.Hyperlinks.Add($ws.Range("C2"), $zhcnXlfUrlOld, "", "", $zhcnOldXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlNew, "", "", $mdNew)
$ws.Hyperlinks.Add($ws.Range("C3"), $zhcnXlfUrlNew, "", "", $zhcnNewXlf)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $config)

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = $mdNew
$ws.Range("B3").Value = $ready
$ws.Range("C3").Value = $dedeNewXlf
$ws.Range("D3").Value = $dedeNewDt
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $include

$ws.Range("A4").Value = $config
$ws.Range("B4").Value = $notLoc
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = $ignored

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlOld, "", "", $mdOld)
$ws.Hyperlinks.Add($ws.Range("C2"), $dedeXlfUrlOld, "", "", $dedeOldXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlNew, "", "", $mdNew)
$ws.Hyperlinks.Add($ws.Range("C3"), $dedeXlfUrlNew, "", "", $dedeNewXlf)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $config)

Write-Output "Report generated for handoff."
